$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refresh publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 ("Contact" / "No display for ContactDetail") becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" row - remove it entirely (rows below shift up)
$meta.Rows.Item(11).Delete()

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (the root Extension element): Slice Name / Label columns (K, L)
# now reflect the resource's own title/description instead of generic text
$elements.Range("K2").Value = "Episode Cost Scale"
$elements.Range("L2").Value = "Cost scale of the episode of care"
